$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price (D) and 1h volume change (E) values
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.398.43'
$ws.Range('E2').Value = '  -3.71%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.954.76'
$ws.Range('E3').Value = '  -1.98%  '
$ws.Range('E4').Value = '  -0.76%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '320.44'
$ws.Range('E5').Value = '  -2.72%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.006'
$ws.Range('E6').Value = '  -0.70%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4752'
$ws.Range('E7').Value = '  -5.41%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4048'
$ws.Range('E8').Value = '  -4.47%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '53.53'
$ws.Range('E9').Value = '  -0.34%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.08394'
$ws.Range('E10').Value = '  -6.18%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.057'
$ws.Range('E11').Value = '  -5.01%  '
$ws.Range('E12').Value = '  -3.90%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.934.09'
$ws.Range('E13').Value = '  -3.17%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.598'
$ws.Range('E14').Value = '  -4.70%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.142'
$ws.Range('E15').Value = '  -4.89%  '
$ws.Range('E16').Value = '  -0.50%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '90.08'
$ws.Range('E17').Value = '  -4.61%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001066'
$ws.Range('E18').Value = '  -4.23%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06583'
$ws.Range('E19').Value = '  -2.68%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.53'
$ws.Range('E20').Value = '  -4.51%  '
$ws.Range('E21').Value = '  -0.67%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.818'
$ws.Range('E22').Value = '  -1.98%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '28.433.62'
$ws.Range('E23').Value = '  -3.72%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.51'
$ws.Range('E24').Value = '  -4.75%  '
$ws.Range('E25').Value = '  -1.07%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.189.22'
$ws.Range('E26').Value = '  -2.25%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '155.18'
$ws.Range('E27').Value = '  -1.27%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.17'
$ws.Range('E28').Value = '  -3.02%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.913'
$ws.Range('E29').Value = '  -6.19%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.151'
$ws.Range('E30').Value = '  -6.69%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '123.45'
$ws.Range('E31').Value = '  -3.25%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.9772'
$ws.Range('E32').Value = '  -7.78%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09594'
$ws.Range('E33').Value = '  -3.50%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.446'
$ws.Range('E34').Value = '  -6.47%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.593'
$ws.Range('E35').Value = '  -4.09%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.655'
$ws.Range('E36').Value = '  -3.60%  '
$ws.Range('E37').Value = '  -3.68%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02328'
$ws.Range('E38').Value = '  -5.47%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06206'
$ws.Range('E39').Value = '  -2.93%  '
$ws.Range('E40').Value = '  -4.18%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6193'
$ws.Range('E41').Value = '  -5.22%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '11.14'
$ws.Range('E42').Value = '  -4.20%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.006'
$ws.Range('E43').Value = '  -0.67%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.1920'
$ws.Range('E44').Value = '  -5.91%  '
$ws.Range('E45').Value = '  +3.52%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5945'
$ws.Range('E46').Value = '  -6.10%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '13.01'
$ws.Range('E47').Value = '  -3.45%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.057'
$ws.Range('E48').Value = '  -6.94%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.387'
$ws.Range('E49').Value = '  -3.31%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.00000000326'
$ws.Range('E50').Value = '  -4.16%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06818'
$ws.Range('E51').Value = '  -2.00%  '
